$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text-number format on columns D and E for the cells we are about to write,
# so numeric-looking strings ("0.999", "3.10", etc.) are stored as text, matching the
# original inlineStr cells, instead of being auto-coerced to Excel numbers.
$dRange = $ws.Range("D2:D51")
$eRange = $ws.Range("E2:E51")
$dRange.NumberFormat = "@"
$eRange.NumberFormat = "@"

$ws.Range('D2').Value = '66.457.06'
$ws.Range('E2').Value = '  -5.97%  '
$ws.Range('D3').Value = '3.195.50'
$ws.Range('E3').Value = '  -9.38%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').Value = '575.31'
$ws.Range('E5').Value = '  -6.34%  '
$ws.Range('D6').Value = '149.19'
$ws.Range('E6').Value = '  -14.20%  '
$ws.Range('D7').Value = '0.998'
$ws.Range('E7').Value = '  -0.11%  '
$ws.Range('D8').Value = '3.190.82'
$ws.Range('E8').Value = '  -9.38%  '
$ws.Range('D9').Value = '0.540'
$ws.Range('E9').Value = '  -11.50%  '
$ws.Range('E10').Value = '  -14.47%  '
$ws.Range('D11').Value = '6.36'
$ws.Range('E11').Value = '  -11.27%  '
$ws.Range('D12').Value = '0.495'
$ws.Range('E12').Value = '  -16.17%  '
$ws.Range('D13').Value = '38.18'
$ws.Range('E13').Value = '  -18.06%  '
$ws.Range('D14').Value = '0.0000240'
$ws.Range('E14').Value = '  -13.48%  '
$ws.Range('D15').Value = '3.693.45'
$ws.Range('E15').Value = '  -9.77%  '
$ws.Range('D16').Value = '66.412.23'
$ws.Range('E16').Value = '  -6.01%  '
$ws.Range('D17').Value = '3.183.17'
$ws.Range('E17').Value = '  -10.10%  '
$ws.Range('E18').Value = '  -6.60%  '
$ws.Range('D19').Value = '529.39'
$ws.Range('E19').Value = '  -14.01%  '
$ws.Range('D20').Value = '7.09'
$ws.Range('E20').Value = '  -16.31%  '
$ws.Range('D21').Value = '14.96'
$ws.Range('E21').Value = '  -15.78%  '
$ws.Range('D22').Value = '0.752'
$ws.Range('E22').Value = '  -15.05%  '
$ws.Range('D23').Value = '7.67'
$ws.Range('E23').Value = '  -15.03%  '
$ws.Range('D24').Value = '84.54'
$ws.Range('E24').Value = '  -14.44%  '
$ws.Range('D25').Value = '13.22'
$ws.Range('E25').Value = '  -15.92%  '
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  -0.04%  '
$ws.Range('D27').Value = '3.10'
$ws.Range('E27').Value = '  -18.11%  '
$ws.Range('E28').Value = '  -17.23%  '
$ws.Range('E29').Value = '  -13.68%  '
$ws.Range('D30').Value = '28.84'
$ws.Range('E30').Value = '  -15.15%  '
$ws.Range('E31').Value = '  -15.73%  '
$ws.Range('D32').Value = '1.12'
$ws.Range('E32').Value = '  -14.05%  '
$ws.Range('B33').Value = 'Bittensor'
$ws.Range('C33').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D33').Value = '529.99'
$ws.Range('E33').Value = '  -14.41%  '
$ws.Range('B34').Value = 'Filecoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D34').Value = '6.46'
$ws.Range('E34').Value = '  -20.91%  '
$ws.Range('D35').Value = '5.61'
$ws.Range('E35').Value = '  -18.22%  '
$ws.Range('E36').Value = '  +0.07%  '
$ws.Range('D37').Value = '52.81'
$ws.Range('E37').Value = '  -7.35%  '
$ws.Range('D38').Value = '0.0851'
$ws.Range('E38').Value = '  -15.29%  '
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').Value = '0.0412'
$ws.Range('E39').Value = '  -16.85%  '
$ws.Range('B40').Value = 'Cosmos'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D40').Value = '9.05'
$ws.Range('E40').Value = '  -16.62%  '
$ws.Range('E41').Value = '  -14.78%  '
$ws.Range('D42').Value = '2.876.97'
$ws.Range('E42').Value = '  -14.84%  '
$ws.Range('E43').Value = '  -25.87%  '
$ws.Range('D44').Value = '0.259'
$ws.Range('E44').Value = '  -17.20%  '
$ws.Range('D45').Value = '0.0₃0579'
$ws.Range('E45').Value = '  -21.88%  '
$ws.Range('D47').Value = '25.76'
$ws.Range('E47').Value = '  -20.24%  '
$ws.Range('E48').Value = '  -21.13%  '
$ws.Range('E49').Value = '  -19.35%  '
$ws.Range('D50').Value = '0.113'
$ws.Range('E50').Value = '  -13.96%  '
$ws.Range('D51').Value = '122.21'
$ws.Range('E51').Value = '  -8.64%  '

# Restore default cell style (removes the temporary text-number format we applied,
# keeping the saved file free of spurious style attributes on these cells).
$dRange.Style = "Normal"
$eRange.Style = "Normal"
